$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1) Update the "总计" (summary) sheet: insert a new row for 2022-Q4
#    at the top of the data (row 2), pushing the existing rows down.
# ---------------------------------------------------------------------
$summary = $wb.Worksheets.Item(1)
$summary.Rows.Item(2).Insert()
$summary.Range("B2:D2").ClearFormats()

# give A2 the same style as the rest of column A (bold / bordered)
$summary.Range("A3").Copy()
$summary.Range("A2").PasteSpecial(-4122)

$summary.Range("A2").Value = 0
$summary.Range("B2").Value = "2022-Q4"
$summary.Range("C2").Value = 36
$summary.Range("D2").Value = 7.13

# column "A" is a running index (0,1,2,...); renumber the rows that
# shifted down one position
$summary.Range("A3").Value = 1
$summary.Range("A4").Value = 2
$summary.Range("A5").Value = 3
$summary.Range("A6").Value = 4

# ---------------------------------------------------------------------
# 2) Create the new "2022-Q4" worksheet by duplicating the existing
#    "2022-Q3" sheet (keeps identical layout / header / styling) and
#    placing the copy immediately before it.
# ---------------------------------------------------------------------
$q3 = $wb.Worksheets.Item(2)
$q3.Copy($q3)
$q4 = $wb.Worksheets.Item(2)
$q4.Name = "2022-Q4"

# extend formatting (column-A bold/border style, etc.) down through row 37
$q4.Rows.Item(15).Copy()
$q4.Range("A16:H37").PasteSpecial(-4122)

$q4Data = @(
    @(0, "012930", "中庚价值先锋股票", "68.71", "94.78", "5.17", "3.5523", 4),
    @(1, "001986", "前海开源人工智能主题灵活配置混合", "7.25", "92.91", "8.78", "0.6366", 4),
    @(2, "005763", "中欧电子信息产业沪港深股票C", "14.81", "91.56", "3.07", "0.4547", 8),
    @(3, "014292", "嘉实产业领先混合A", "13.07", "92.57", "3.29", "0.4300", 6),
    @(4, "012985", "平安优势回报1年持有混合A", "12.27", "94.40", "2.99", "0.3669", 9),
    @(5, "004616", "中欧电子信息产业沪港深股票A", "6.80", "91.56", "3.07", "0.2088", 8),
    @(6, "012917", "平安优势领航1年持有期混合A", "6.72", "93.97", "2.98", "0.2003", 9),
    @(7, "002450", "平安睿享文娱灵活配置混合A", "3.85", "94.03", "4.55", "0.1752", 3),
    @(8, "001103", "前海开源工业革命4.0灵活配置混合", "4.12", "87.14", "3.32", "0.1368", 8),
    @(9, "010126", "平安价值成长混合A", "3.21", "94.03", "2.82", "0.0905", 10),
    @(10, "002451", "平安睿享文娱灵活配置混合C", "1.97", "94.03", "4.55", "0.0896", 3),
    @(11, "006101", "平安优势产业灵活配置混合C", "2.85", "94.90", "3.02", "0.0861", 10),
    @(12, "580001", "东吴嘉禾优势精选混合A", "2.06", "87.72", "3.79", "0.0781", 10),
    @(13, "004784", "招商稳健优选股票", "3.60", "90.24", "2.14", "0.0770", 6),
    @(14, "011828", "平安睿享成长混合A", "2.42", "92.96", "2.84", "0.0687", 10),
    @(15, "501099", "平安科技创新 3 年封闭混合", "2.73", "91.32", "2.46", "0.0672", 9),
    @(16, "006100", "平安优势产业灵活配置混合A", "2.15", "94.90", "3.02", "0.0649", 10),
    @(17, "009008", "平安科技创新混合A", "2.49", "92.02", "2.46", "0.0613", 9),
    @(18, "013687", "平安成长龙头1年持有混合A", "1.21", "94.87", "3.53", "0.0427", 10),
    @(19, "010127", "平安价值成长混合C", "1.33", "94.03", "2.82", "0.0375", 10),
    @(20, "012986", "平安优势回报1年持有混合C", "1.14", "94.40", "2.99", "0.0341", 9),
    @(21, "011829", "平安睿享成长混合C", "1.09", "92.96", "2.84", "0.0310", 10),
    @(22, "009009", "平安科技创新混合C", "0.95", "92.02", "2.46", "0.0234", 9),
    @(23, "007518", "东方阿尔法优选混合A", "0.94", "92.44", "2.09", "0.0196", 5),
    @(24, "013688", "平安成长龙头1年持有混合C", "0.51", "94.87", "3.53", "0.0180", 10),
    @(25, "004321", "前海开源沪港深强国产业灵活配置混合", "0.33", "90.19", "4.32", "0.0143", 8),
    @(26, "007894", "平安估值精选混合C", "0.50", "94.39", "2.85", "0.0142", 10),
    @(27, "000679", "招商丰利灵活配置混合A", "0.29", "84.08", "3.74", "0.0108", 8),
    @(28, "007893", "平安估值精选混合A", "0.34", "94.39", "2.85", "0.0097", 10),
    @(29, "007519", "东方阿尔法优选混合C", "0.41", "92.44", "2.09", "0.0086", 5),
    @(30, "014293", "嘉实产业领先混合C", "0.26", "92.57", "3.29", "0.0086", 6),
    @(31, "700004", "平安灵活配置混合A", "0.32", "79.07", "2.44", "0.0078", 8),
    @(32, "015078", "平安灵活配置混合C", "0.18", "79.07", "2.44", "0.0044", 8),
    @(33, "012918", "平安优势领航1年持有期混合C", "0.08", "93.97", "2.98", "0.0024", 9),
    @(34, "002416", "招商丰利灵活配置混合C", "0.02", "84.08", "3.74", "0.0007", 8),
    @(35, "015152", "东吴嘉禾优势精选混合C", "0.01", "87.72", "3.79", "0.0004", 10)
)

# force columns B:G to be treated as text so numeric-looking values
# (e.g. "68.71", "012930") keep their original formatting / leading zeros
$q4.Range("B2:G37").NumberFormat = "@"

for ($i = 0; $i -lt $q4Data.Count; $i++) {
    $row = $q4Data[$i]
    $r = $i + 2

    $q4.Cells.Item($r, 1).Value = $row[0]
    $q4.Cells.Item($r, 2).Value = $row[1]
    $q4.Cells.Item($r, 3).Value = $row[2]
    $q4.Cells.Item($r, 4).Value = $row[3]
    $q4.Cells.Item($r, 5).Value = $row[4]
    $q4.Cells.Item($r, 6).Value = $row[5]
    $q4.Cells.Item($r, 7).Value = $row[6]
    $q4.Cells.Item($r, 8).Value = $row[7]
}

# drop the temporary text number-format now that the values are set, so
# the cells end up with the same (default) style as the source sheet
$q4.Range("B2:G37").ClearFormats()

Write-Host "2022-Q4 sheet populated"
